$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87, shifting existing rows 87-159 down to 88-160.
$ws.Rows("87").Insert()

# Populate the newly inserted row 87 with the new data record.
$ws.Range("A87").Value = 8
$ws.Range("B87").Value = "Terminal La Palmera de La Serena"
$ws.Range("C87").Value = "Coquimbo"
$ws.Range("D87").Value = 44977
$ws.Range("E87").Value = 4
$ws.Range("F87").Value = 100112052
$ws.Range("G87").Value = "Albahaca"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 800
$ws.Range("K87").Value = 4800
$ws.Range("L87").Value = 5000
$ws.Range("M87").Value = 4900
$ws.Range("N87").Value = "$/docena de matas"
$ws.Range("O87").Value = "Provincia del Elquí"
$ws.Range("P87").Value = 817
$ws.Range("Q87").Value = 6
$ws.Range("R87").Value = "Hortaliza"
